$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading-percent values for the "Case with 380 kV" scenario (rows 2-25, columns B,C,D,E,G,L,O)
$ws.Range("B2").Value = 18.57454011872788
$ws.Range("C2").Value = 9.853855477746013
$ws.Range("D2").Value = 5.996619807621492
$ws.Range("E2").Value = 11.26074743519876
$ws.Range("G2").Value = 3.65483739086045
$ws.Range("L2").Value = 9.917129744464303
$ws.Range("O2").Value = 25.81311950481265
$ws.Range("B3").Value = 18.00325837333569
$ws.Range("C3").Value = 9.566954487246271
$ws.Range("D3").Value = 5.880471611436008
$ws.Range("E3").Value = 11.30748871448804
$ws.Range("G3").Value = 3.657698605693289
$ws.Range("L3").Value = 9.890527293875826
$ws.Range("O3").Value = 25.87901689110492
$ws.Range("B4").Value = 17.64595112973065
$ws.Range("C4").Value = 9.385235561588056
$ws.Range("D4").Value = 5.809853088783401
$ws.Range("E4").Value = 11.33816515830682
$ws.Range("G4").Value = 3.659547268771552
$ws.Range("L4").Value = 9.875970262351704
$ws.Range("O4").Value = 25.9277208783482
$ws.Range("B5").Value = 17.49892536856162
$ws.Range("C5").Value = 9.30986173853268
$ws.Range("D5").Value = 5.781294205065535
$ws.Range("E5").Value = 11.35116321038829
$ws.Range("G5").Value = 3.660323794785724
$ws.Range("L5").Value = 9.87048885943155
$ws.Range("O5").Value = 25.94962969260857
$ws.Range("B6").Value = 17.47443293795263
$ws.Range("C6").Value = 9.297268491610579
$ws.Range("D6").Value = 5.776566509534832
$ws.Range("E6").Value = 11.35335155388083
$ws.Range("G6").Value = 3.66045413876695
$ws.Range("L6").Value = 9.869606003259277
$ws.Range("O6").Value = 25.95339180583249
$ws.Range("B7").Value = 17.64397373137707
$ws.Range("C7").Value = 9.384224290483145
$ws.Range("D7").Value = 5.809466990418624
$ws.Range("E7").Value = 11.33833844163613
$ws.Range("G7").Value = 3.659557647308134
$ws.Range("L7").Value = 9.875894508535241
$ws.Range("O7").Value = 25.92800801655958
$ws.Range("B8").Value = 18.37904836088398
$ws.Range("C8").Value = 9.756129934898757
$ws.Range("D8").Value = 5.956453261936004
$ws.Range("E8").Value = 11.27645347296665
$ws.Range("G8").Value = 3.655804919548494
$ws.Range("L8").Value = 9.90759084695903
$ws.Range("O8").Value = 25.83412366390773
$ws.Range("B9").Value = 19.75965896433191
$ws.Range("C9").Value = 10.43828642078554
$ws.Range("D9").Value = 6.248334215754381
$ws.Range("E9").Value = 11.17078503827269
$ws.Range("G9").Value = 3.649171072092589
$ws.Range("L9").Value = 9.98366000046599
$ws.Range("O9").Value = 25.71588331151099
$ws.Range("B10").Value = 20.72604163555635
$ws.Range("C10").Value = 10.907110679823
$ws.Range("D10").Value = 6.462532023255845
$ws.Range("E10").Value = 11.10271458485408
$ws.Range("G10").Value = 3.644734168103836
$ws.Range("L10").Value = 10.04775546600764
$ws.Range("O10").Value = 25.66974808539241
$ws.Range("B11").Value = 21.15330664090591
$ws.Range("C11").Value = 11.11272107737919
$ws.Range("D11").Value = 6.559426670798619
$ws.Range("E11").Value = 11.07382481894056
$ws.Range("G11").Value = 3.642809495978688
$ws.Range("L11").Value = 10.07863008999898
$ws.Range("O11").Value = 25.6577111466097
$ws.Range("B12").Value = 21.31318410860069
$ws.Range("C12").Value = 11.18943248093371
$ws.Range("D12").Value = 6.596000849106912
$ws.Range("E12").Value = 11.06318365119712
$ws.Range("G12").Value = 3.642094062673491
$ws.Range("L12").Value = 10.0905622549998
$ws.Range("O12").Value = 25.65444716646762
$ws.Range("B13").Value = 21.27883907642959
$ws.Range("C13").Value = 11.17296309838901
$ws.Range("D13").Value = 6.588129785810466
$ws.Range("E13").Value = 11.06546212533001
$ws.Range("G13").Value = 3.642247549341065
$ws.Range("L13").Value = 10.08798185266379
$ws.Range("O13").Value = 25.65509246328047
$ws.Range("B14").Value = 21.16649907101365
$ws.Range("C14").Value = 11.11905544663543
$ws.Range("D14").Value = 6.562438227879348
$ws.Range("E14").Value = 11.07294337452906
$ws.Range("G14").Value = 3.642750368759569
$ws.Range("L14").Value = 10.07960697154528
$ws.Range("O14").Value = 25.65741664332864
$ws.Range("B15").Value = 21.09743374424196
$ws.Range("C15").Value = 11.0858845321524
$ws.Range("D15").Value = 6.546684943016118
$ws.Range("E15").Value = 11.0775647710024
$ws.Range("G15").Value = 3.643060102791474
$ws.Range("L15").Value = 10.07450825669059
$ws.Range("O15").Value = 25.65900900022854
$ws.Range("B16").Value = 20.69785782774168
$ws.Range("C16").Value = 10.89351525508572
$ws.Range("D16").Value = 6.456185353227381
$ws.Range("E16").Value = 11.10464439156339
$ws.Range("G16").Value = 3.644861828939436
$ws.Range("L16").Value = 10.04577175958613
$ws.Range("O16").Value = 25.67071546106969
$ws.Range("B17").Value = 20.44946183790688
$ws.Range("C17").Value = 10.77350568803923
$ws.Range("D17").Value = 6.400498562259978
$ws.Range("E17").Value = 11.12178870807851
$ws.Range("G17").Value = 3.645991074052822
$ws.Range("L17").Value = 10.02857829173243
$ws.Range("O17").Value = 25.68019475382157
$ws.Range("B18").Value = 20.30543638379748
$ws.Range("C18").Value = 10.70376083537996
$ws.Range("D18").Value = 6.36841992540591
$ws.Range("E18").Value = 11.13184503525816
$ws.Range("G18").Value = 3.646649409725344
$ws.Range("L18").Value = 10.01885110275808
$ws.Range("O18").Value = 25.6864890672003
$ws.Range("B19").Value = 20.25647813779662
$ws.Range("C19").Value = 10.6800245459305
$ws.Range("D19").Value = 6.35755145472502
$ws.Range("E19").Value = 11.13528348224052
$ws.Range("G19").Value = 3.646873828664495
$ws.Range("L19").Value = 10.01558566125242
$ws.Range("O19").Value = 25.68876460685444
$ws.Range("B20").Value = 20.47602460164856
$ws.Range("C20").Value = 10.78635563485001
$ws.Range("D20").Value = 6.406431886044701
$ws.Range("E20").Value = 11.11994344446209
$ws.Range("G20").Value = 3.645869951304044
$ws.Range("L20").Value = 10.03039184006936
$ws.Range("O20").Value = 25.67909846036357
$ws.Range("B21").Value = 21.19954918551632
$ws.Range("C21").Value = 11.13492095971766
$ws.Range("D21").Value = 6.569987959583202
$ws.Range("E21").Value = 11.0707378398854
$ws.Range("G21").Value = 3.642602315481541
$ws.Range("L21").Value = 10.08206039900577
$ws.Range("O21").Value = 25.6566988010034
$ws.Range("B22").Value = 21.66117601289356
$ws.Range("C22").Value = 11.3560141573681
$ws.Range("D22").Value = 6.676178359574497
$ws.Range("E22").Value = 11.04032066613174
$ws.Range("G22").Value = 3.640544787211315
$ws.Range("L22").Value = 10.11722811045998
$ws.Range("O22").Value = 25.64960516207883
$ws.Range("B23").Value = 21.41586872389578
$ws.Range("C23").Value = 11.238641223035
$ws.Range("D23").Value = 6.61957932802222
$ws.Range("E23").Value = 11.05639544339857
$ws.Range("G23").Value = 3.641635811007331
$ws.Range("L23").Value = 10.09833257881946
$ws.Range("O23").Value = 25.65269864827907
$ws.Range("B24").Value = 20.46401936294941
$ws.Range("C24").Value = 10.78054850565129
$ws.Range("D24").Value = 6.403749625357507
$ws.Range("E24").Value = 11.12077706569447
$ws.Range("G24").Value = 3.645924682490954
$ws.Range("L24").Value = 10.02957144359834
$ws.Range("O24").Value = 25.67959146436484
$ws.Range("B25").Value = 19.39386940208042
$ws.Range("C25").Value = 10.25918930894648
$ws.Range("D25").Value = 6.169242305875739
$ws.Range("E25").Value = 11.19769168827295
$ws.Range("G25").Value = 3.650888594839309
$ws.Range("L25").Value = 9.961618575114384
$ws.Range("O25").Value = 25.74075538013524

Write-Host "Updated 168 cells across rows 2-25 for columns B,C,D,E,G,L,O"
